$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1999
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1999
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5997
$ws.Range("N39").Value = -6589
$ws.Range("M39").Value = $null

$ws.Range("H40").Value = 3750
$ws.Range("I40").Value = 2500.25
$ws.Range("J40").Value = 4999.75
$ws.Range("K40").Value = 2500.25
$ws.Range("L40").Value = 4999.75
$ws.Range("M40").Value = -2325.25
$ws.Range("N40").Value = -5349.75

$ws.Range("H41").Value = 407.33334
$ws.Range("I41").Value = 103.2
$ws.Range("K41").Value = 103.2
$ws.Range("M41").Value = 336.8

$ws.Range("H53").Value = 653
$ws.Range("I53").Value = 538
$ws.Range("K53").Value = 538
$ws.Range("M53").Value = 99

$ws.Range("H74").Value = 3713.8333
$ws.Range("I74").Value = 3713.8333
$ws.Range("K74").Value = 3713.8333
$ws.Range("M74").Value = -2777.8333

$ws.Range("H77").Value = 3713.8333
$ws.Range("I77").Value = 3713.8333
$ws.Range("K77").Value = 18569.1665
$ws.Range("M77").Value = -13889.1665

$ws.Range("H105").Value = 58571.145
$ws.Range("J105").Value = 59999.668
$ws.Range("L105").Value = 59999.668
$ws.Range("N105").Value = -66987.66800000001

$ws.Range("H107").Value = 855138.3
$ws.Range("I107").Value = 926233.2
$ws.Range("K107").Value = 926233.2
$ws.Range("M107").Value = -924313.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 8750869
$ws.Range("J13").Value = 1158.1666
$ws.Range("L13").Value = 1158.1666
$ws.Range("N13").Value = -1446.1666

$ws.Range("H32").Value = 2947.2917
$ws.Range("I32").Value = 2829.068
$ws.Range("K32").Value = 2829.068
$ws.Range("M32").Value = -2542.068

$ws.Range("H61").Value = 6225.6875
$ws.Range("I61").Value = 6484.385
$ws.Range("K61").Value = 6484.385
$ws.Range("M61").Value = -6272.385

$ws.Range("H97").Value = 656.8570999999999
$ws.Range("I97").Value = 595.2
$ws.Range("K97").Value = 595.2
$ws.Range("M97").Value = -99.20000000000005

$ws.Range("H103").Value = 30329
$ws.Range("I103").Value = 30329
$ws.Range("K103").Value = 30329
$ws.Range("M103").Value = -29157

$ws.Range("H110").Value = 1015.38464
$ws.Range("I110").Value = 1015.38464
$ws.Range("K110").Value = 1015.38464
$ws.Range("M110").Value = 1029.61536

$ws.Range("H136").Value = 6225.6875
$ws.Range("I136").Value = 6484.385
$ws.Range("K136").Value = 19453.155
$ws.Range("M136").Value = -16903.155

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3003.6765
$ws.Range("I94").Value = 2586.7856
$ws.Range("K94").Value = 2586.7856
$ws.Range("M94").Value = -2135.7856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 134.6
$ws.Range("I7").Value = 57.81818
$ws.Range("J7").Value = 228.44444
$ws.Range("K7").Value = 57.81818
$ws.Range("L7").Value = 228.44444
$ws.Range("M7").Value = 55.18182
$ws.Range("N7").Value = -454.44444

$ws.Range("H31").Value = 4434.7036
$ws.Range("I31").Value = 1756.3684
$ws.Range("K31").Value = 1756.3684
$ws.Range("M31").Value = -1461.3684

$ws.Range("H34").Value = 4434.7036
$ws.Range("I34").Value = 1756.3684
$ws.Range("K34").Value = 1756.3684
$ws.Range("M34").Value = -1554.3684

$ws.Range("H42").Value = 10033.333
$ws.Range("I42").Value = 10033.333
$ws.Range("K42").Value = 10033.333
$ws.Range("M42").Value = -9440.333000000001

$ws.Range("H122").Value = 3308.0334
$ws.Range("I122").Value = 925.8
$ws.Range("J122").Value = 4499.15
$ws.Range("K122").Value = 2777.4
$ws.Range("L122").Value = 13497.45
$ws.Range("M122").Value = -327.3999999999996
$ws.Range("N122").Value = -18397.45

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 91.19231000000001
$ws.Range("J12").Value = 107.72222
$ws.Range("L12").Value = 323.16666
$ws.Range("N12").Value = -669.16666

$ws.Range("H13").Value = 1222
$ws.Range("J13").Value = 3500
$ws.Range("L13").Value = 10500
$ws.Range("N13").Value = -10836

$ws.Range("H131").Value = 2293.1875
$ws.Range("J131").Value = 2428.3447
$ws.Range("L131").Value = 7285.034100000001
$ws.Range("N131").Value = -17365.0341

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6339.0625
$ws.Range("I122").Value = 5998.3105
$ws.Range("K122").Value = 17994.9315
$ws.Range("M122").Value = -15544.9315

$ws.Range("H123").Value = 22499.125
$ws.Range("J123").Value = 22499.125
$ws.Range("L123").Value = 22499.125
$ws.Range("N123").Value = -27399.125

$ws.Range("H126").Value = 2485.3125
$ws.Range("J126").Value = 3059.1428
$ws.Range("L126").Value = 9177.428400000001
$ws.Range("N126").Value = -14117.4284

$ws.Range("H132").Value = 2404.8823
$ws.Range("J132").Value = 1561
$ws.Range("L132").Value = 4683
$ws.Range("N132").Value = -9743

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1233.1666
$ws.Range("I9").Value = 1833
$ws.Range("J9").Value = 633.3333
$ws.Range("K9").Value = 1833
$ws.Range("L9").Value = 633.3333
$ws.Range("M9").Value = -1609
$ws.Range("N9").Value = -1081.3333

$ws.Range("H10").Value = 2922.3845
$ws.Range("I10").Value = 152.66667
$ws.Range("J10").Value = 3753.3
$ws.Range("K10").Value = 152.66667
$ws.Range("L10").Value = 3753.3
$ws.Range("M10").Value = -12.66667000000001
$ws.Range("N10").Value = -4033.3

$ws.Range("H11").Value = 7497
$ws.Range("J11").Value = 7497
$ws.Range("L11").Value = 7497
$ws.Range("N11").Value = -7777

$ws.Range("H12").Value = 2593.875
$ws.Range("I12").Value = 347.66666
$ws.Range("J12").Value = 3941.6
$ws.Range("K12").Value = 347.66666
$ws.Range("L12").Value = 3941.6
$ws.Range("N12").Value = -4281.6
$ws.Range("M12").Value = -177.66666

$ws.Range("H14").Value = 11333
$ws.Range("J14").Value = 20000
$ws.Range("L14").Value = 20000
$ws.Range("N14").Value = -20344

$ws.Range("H17").Value = 1324.9744
$ws.Range("I17").Value = 1250.35
$ws.Range("J17").Value = 1403.5264
$ws.Range("K17").Value = 1250.35
$ws.Range("L17").Value = 1403.5264
$ws.Range("M17").Value = -1080.35
$ws.Range("N17").Value = -1743.5264

$ws.Range("H46").Value = 843.2857
$ws.Range("I46").Value = 700
$ws.Range("K46").Value = 700
$ws.Range("M46").Value = -512

$ws.Range("H93").Value = 3974.3215
$ws.Range("I93").Value = 2207.1538
$ws.Range("J93").Value = 5505.8667
$ws.Range("K93").Value = 2207.1538
$ws.Range("L93").Value = 5505.8667
$ws.Range("M93").Value = -959.1538
$ws.Range("N93").Value = -8001.8667

$ws.Range("H100").Value = 12502
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 314.33334
$ws.Range("I100").Value = 243.25
$ws.Range("K100").Value = 486.5
$ws.Range("M100").Value = 54.5
